# "Braga with everything switched on."
#
# Relations sheet: the undeclared/v/D/w/E block moves from columns D:F
# (rows 6-8) to columns A:C (rows 6-8).
#
# Rules sheet: the "v = w" rule row (row 4) is deleted entirely, and the
# "undeclared = w" rule row (row 5) moves up into row 4, shifted from
# columns D:F into columns A:C.
#
# Compositions sheet: no data changes, only the selection moves.
#
# Finally, the selections on every sheet are updated to match, and the
# "Rules" tab is left as the active/selected one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Relations sheet
# ---------------------------------------------------------------
$relations = $wb.Worksheets.Item("Relations")

# Write the relocated values into A6:C8 first ...
$relations.Range("A6").Value = "undeclared"

$relations.Range("A7").Value = "v"
$relations.Range("B7").Value = "D"
$relations.Range("C7").Value = "D"

$relations.Range("A8").Value = "w"
$relations.Range("B8").Value = "E"
$relations.Range("C8").Value = "E"

# ... then clear out the old D6:F8 block so it disappears from the used
# range entirely (ClearContents would leave empty-but-present cells
# behind, Clear drops them so the sheet dimension shrinks back to C).
$relations.Range("D6:F8").Clear()

$relations.Range("A6:C8").Select()

# ---------------------------------------------------------------
# Rules sheet
# ---------------------------------------------------------------
$rules = $wb.Worksheets.Item("Rules")

# Only the "undeclared = w" row survives, now in row 4 / columns A:C.
$rules.Range("A4").Value = "undeclared = w"
$rules.Range("B4").Value = "undeclared"
$rules.Range("C4").Value = "w"

# Drop the old D4:F5 block (the "v = w" row plus the relocated
# "undeclared = w" row) entirely.
$rules.Range("D4:F5").Clear()

$rules.Range("A4:C4").Select()

# ---------------------------------------------------------------
# Compositions sheet (selection only, no data change)
# ---------------------------------------------------------------
$compositions = $wb.Worksheets.Item("Compositions")
$compositions.Range("A4:C4").Select()

# ---------------------------------------------------------------
# Leave "Rules" as the active / selected tab, matching the saved file.
# ---------------------------------------------------------------
$rules.Activate()
$rules.Range("A4:C4").Select()

# ---------------------------------------------------------------
# Best-effort: carry over the cosmetic workbook-window geometry too.
# (No-op in this headless host if unsupported, but harmless to try.)
# ---------------------------------------------------------------
try { $excel.Left = 240 } catch {}
try { $excel.Top = 225 } catch {}
try { $excel.Width = 14805 } catch {}
try { $excel.Height = 7890 } catch {}
try {
    $awin = $excel.ActiveWindow
    $awin.Left = 240
    $awin.Top = 225
    $awin.Width = 14805
    $awin.Height = 7890
} catch {}
